# Updated cryptos list on Mon Feb 12 17:59:32 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $text) {
    # Force the cell to be stored as text (matches original inlineStr cells),
    # even when the text looks like a number (e.g. "110.07").
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "49.867.02"
$ws.Range("E2").Value = "  +3.29%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "2.606.94"
$ws.Range("E3").Value = "  +3.86%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - was BNB, now Solana
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws "D5" "110.07"
$ws.Range("E5").Value = "  +1.21%  "

# Row 6 - was Solana, now BNB
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue $ws "D6" "323.63"
$ws.Range("E6").Value = "  +0.69%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.79%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Cardano
Set-TextValue $ws "D9" "0.564"
$ws.Range("E9").Value = "  +3.57%  "

# Row 10 - Avalanche
Set-TextValue $ws "D10" "40.87"
$ws.Range("E10").Value = "  +2.27%  "

# Row 11 - Chainlink
Set-TextValue $ws "D11" "20.73"
$ws.Range("E11").Value = "  +3.21%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  +0.47%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.55%  "

# Row 14 - Polkadot
Set-TextValue $ws "D14" "7.37"
$ws.Range("E14").Value = "  +2.24%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D15" "3.016.85"
$ws.Range("E15").Value = "  +3.81%  "

# Row 16 - WrappedEther
Set-TextValue $ws "D16" "2.626.42"
$ws.Range("E16").Value = "  +4.58%  "

# Row 18 - WrappedBTC
Set-TextValue $ws "D18" "49.819.35"
$ws.Range("E18").Value = "  +3.52%  "

# Row 19 - ImmutableX
$ws.Range("E19").Value = "  +11.63%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("E20").Value = "  +1.90%  "

# Row 21 - Uniswap
Set-TextValue $ws "D21" "6.78"
$ws.Range("E21").Value = "  +0.73%  "

# Row 22 - ShibaInu
Set-TextValue $ws "D22" "0.0₃0951"
$ws.Range("E22").Value = "  +0.16%  "

# Row 23 - BitcoinCash
Set-TextValue $ws "D23" "283.65"
$ws.Range("E23").Value = "  +2.14%  "

# Row 24 - Litecoin
Set-TextValue $ws "D24" "72.85"
$ws.Range("E24").Value = "  +0.91%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +0.33%  "

# Row 26 - EthereumClassic
Set-TextValue $ws "D26" "26.65"
$ws.Range("E26").Value = "  +3.11%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.00%  "

# Row 28 - Kaspa
Set-TextValue $ws "D28" "0.146"
$ws.Range("E28").Value = "  +4.04%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -7.02%  "

# Row 30 - Cosmos
$ws.Range("E30").Value = "  +1.40%  "

# Row 31 - InjectiveProtocol
$ws.Range("E31").Value = "  +0.84%  "

# Row 32 - OKB
Set-TextValue $ws "D32" "49.54"
$ws.Range("E32").Value = "  +0.70%  "

# Row 33 - Celestia
Set-TextValue $ws "D33" "19.67"
$ws.Range("E33").Value = "  +0.45%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  +1.51%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  -0.17%  "

# Row 36 - Hedera
$ws.Range("E36").Value = "  +0.94%  "

# Row 37 - ARBITRUM
Set-TextValue $ws "D37" "2.05"
$ws.Range("E37").Value = "  +4.64%  "

# Row 38 - RenderToken
Set-TextValue $ws "D38" "4.74"
$ws.Range("E38").Value = "  +1.89%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +3.80%  "

# Row 40 - EnergySwap
$ws.Range("E40").Value = "  +6.45%  "

# Row 41 - Monero
Set-TextValue $ws "D41" "124.68"
$ws.Range("E41").Value = "  +2.01%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  +0.68%  "

# Row 43 - WEMIXToken
$ws.Range("E43").Value = "  +0.32%  "

# Row 44 - VeChain
Set-TextValue $ws "D44" "0.0316"
$ws.Range("E44").Value = "  +2.72%  "

# Row 45 - NEARProtocol
$ws.Range("E45").Value = "  +5.23%  "

# Row 46 - Maker
Set-TextValue $ws "D46" "2.041.78"
$ws.Range("E46").Value = "  +1.98%  "

# Row 47 - Stacks
Set-TextValue $ws "D47" "2.02"
$ws.Range("E47").Value = "  +8.81%  "

# Row 48 - ApeXProtocol
$ws.Range("E48").Value = "  +8.60%  "

# Row 49 - FraxShare
Set-TextValue $ws "D49" "9.11"
$ws.Range("E49").Value = "  +1.11%  "

# Row 50 - THORChain
$ws.Range("E50").Value = "  +2.84%  "

# Row 51 - BitcoinSV
Set-TextValue $ws "D51" "81.64"
$ws.Range("E51").Value = "  +1.65%  "
